$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Показники"
$ws.Range("B1").Value = "Попередній період"
$ws.Range("C1").Value = "Поточний період"

# Data rows
$data = @(
    @("Дохід від реалізації", 200, 300),
    @("Собівартість", 150, 170),
    @("Валовий", 50, 130),
    @("Адміністративні витрати", 50, 60),
    @("Витрати на збут", 40, 30),
    @("Інші витрати", 80, 100),
    @("Інші доходи", 60, 80),
    @("Інші фінансові доходи", 120, 200),
    @("Фінансовий результат", 60, 220)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Apply header style (bold, bordered, centered) to B1:C1 to match A1's existing style
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)

$ws.Range("G18").Select()
